$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 98: correct date/time (A98) and the high price (C98) ---
$ws.Cells.Item(98, 1).Value = 45483.2916666667
$ws.Cells.Item(98, 3).Value = 3.33999991416931

# --- Append the newest observation as row 99 ---
$ws.Cells.Item(99, 2).Value = 8250
$ws.Cells.Item(99, 3).Value = 3.25
$ws.Cells.Item(99, 4).Value = 3.13000011444092
$ws.Cells.Item(99, 5).Value = 3.22000002861023
$ws.Cells.Item(99, 6).Value = 3.24000000953674
$ws.Cells.Item(99, 8).Value = "ESPE.MI"

# adj_close (column G) is stored as text for every row in this sheet, so force
# a text/shared-string cell instead of a number, then re-copy the plain
# (unformatted) style from the cell above so no stray number format lingers.
$ws.Cells.Item(99, 7).NumberFormat = "@"
$ws.Cells.Item(99, 7).Value = "3.24000000953674"
$ws.Cells.Item(98, 7).Copy()
$ws.Cells.Item(99, 7).PasteSpecial(-4122)

# Date column (A) uses a custom date/time display format - copy that
# formatting (not just the value) from the row above so A99 matches the rest
# of the column.
$ws.Cells.Item(99, 1).Value = 45484.6250347222
$ws.Cells.Item(98, 1).Copy()
$ws.Cells.Item(99, 1).PasteSpecial(-4122)
